# Bulk order template update: new shared-string labels (ReferenceNumber,
# FromPhone/ToPhone wording, service-speed descriptions), reorganised
# header row (Company columns dropped, ReferenceNumber added), wrap-text
# styling on the header rows, resized columns, and a rebuilt "Service
# Speeds" lookup sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet 1 ("kemlabels-bulk-order-template")
# ---------------------------------------------------------------------

# Row 1: courier / service-speed / signature picks
$ws1.Range("A1").Value = "USPS"
$ws1.Range("B1").Value = "Ground Advantage: 1-5 days"
$ws1.Range("C1").Value = "YES"

# Row 2: column headers (FromCompany / ToCompany dropped, phone fields
# renamed + repositioned, ReferenceNumber columns added at the end,
# Description2 dropped)
$ws1.Range("A2").Value = "FromCountry"
$ws1.Range("B2").Value = "FromName"
$ws1.Range("C2").Value = "FromPhone (Optional)"
$ws1.Range("D2").Value = "FromStreet1"
$ws1.Range("E2").Value = "FromStreet2 (Optional)"
$ws1.Range("F2").Value = "FromCity"
$ws1.Range("G2").Value = "FromZip"
$ws1.Range("H2").Value = "FromState"
$ws1.Range("I2").Value = "ToCountry"
$ws1.Range("J2").Value = "ToName"
$ws1.Range("K2").Value = "ToPhone (Optional)"
$ws1.Range("L2").Value = "ToStreet1"
$ws1.Range("M2").Value = "ToStreet2 (Optional)"
$ws1.Range("N2").Value = "ToCity"
$ws1.Range("O2").Value = "ToZip"
$ws1.Range("P2").Value = "ToState"
$ws1.Range("Q2").Value = "Length"
$ws1.Range("R2").Value = "Height"
$ws1.Range("S2").Value = "Width"
$ws1.Range("T2").Value = "Weight"
$ws1.Range("U2").Value = "Description (Optional)"
$ws1.Range("V2").Value = "ReferenceNumber"
$ws1.Range("W2").Value = "ReferenceNumber (Optional)"

# The sheet used to run out to column X - drop the now-unused trailing
# column(s) so the used range / dimension shrinks back down to W.
$ws1.Range("X1:X2").ClearContents()

# Explicit row height on the header row (matches the authored file).
$ws1.Rows.Item(2).RowHeight = 15

# Wrap text across both header rows (keeps C1's existing "left align").
$ws1.Range("A1:C1").WrapText = $true
$ws1.Range("A2:W2").WrapText = $true

# Column widths (character units) - the COM width model snaps to the
# sheet's pixel grid, so these are the closest values Excel itself would
# store for the authored widths.
$ws1.Columns.Item(1).ColumnWidth  = 11.6666666666667
$ws1.Columns.Item(2).ColumnWidth  = 27.8333333333333
$ws1.Columns.Item(3).ColumnWidth  = 20.1666666666667
$ws1.Columns.Item(4).ColumnWidth  = 13
$ws1.Columns.Item(5).ColumnWidth  = 20.6666666666667
$ws1.Columns.Item(6).ColumnWidth  = 12.6666666666667
$ws1.Columns.Item(7).ColumnWidth  = 9.5
$ws1.Columns.Item(8).ColumnWidth  = 9.16666666666667
$ws1.Columns.Item(9).ColumnWidth  = 9.33333333333333
$ws1.Columns.Item(10).ColumnWidth = 19.1666666666667
$ws1.Columns.Item(11).ColumnWidth = 18.3333333333333
$ws1.Columns.Item(12).ColumnWidth = 17.6666666666667
$ws1.Columns.Item(13).ColumnWidth = 18.6666666666667
$ws1.Columns.Item(14).ColumnWidth = 18.1666666666667
$ws1.Columns.Item(15).ColumnWidth = 5.66666666666667
$ws1.Columns.Item(16).ColumnWidth = 7.66666666666667
$ws1.Columns.Item(17).ColumnWidth = 6.66666666666667
$ws1.Columns.Item(18).ColumnWidth = 6.66666666666667
$ws1.Columns.Item(19).ColumnWidth = 6.16666666666667
$ws1.Columns.Item(20).ColumnWidth = 6.83333333333333
$ws1.Columns.Item(21).ColumnWidth = 22.6666666666667
$ws1.Columns.Item(22).ColumnWidth = 20.1666666666667
$ws1.Columns.Item(23).ColumnWidth = 26.8333333333333

# Drop the second "YES/NO"-style helper validation rule that pointed at
# the old $A$6:$A$8 seed range - it's gone from the rebuilt sheet.
$ws1.Range("E6").Validation.Delete()

# Portrait page orientation was turned on for printing.
$ws1.PageSetup.Orientation = 1

# Restore the (single-area) active selection.
$ws1.Range("D6,D9").Select()

# ---------------------------------------------------------------------
# Sheet 2 ("Service Speeds") - rebuilt with richer, per-courier speed
# descriptions instead of plain service names.
# ---------------------------------------------------------------------

$ws2.Range("A1").Value = "Ground Advantage: 1-5 days"
$ws2.Range("B1").Value = "Next Day Air Early: 1 day"
$ws2.Range("C1").Value = "Express Early: 1 day"

$ws2.Range("A2").Value = "Priority: 1-3 days"
$ws2.Range("B2").Value = "Next Day Air: 1 day"
$ws2.Range("C2").Value = "Express: 1 day"

$ws2.Range("A3").Value = "Express: 1-2 days"
$ws2.Range("B3").Value = "2nd Day Air: 2 days"
$ws2.Range("C3").Value = "Express Saver: 1 day"

$ws2.Range("B4").Value = "3 Day Select: 3 days"
$ws2.Range("C4").Value = "Expedited: 2 days"

$ws2.Range("B5").Value = "Ground: Min 3 days"
$ws2.Range("C5").Value = "Standard: Flexible"

$ws2.Range("C7").Select()
